$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 8-15 ---
# Row 8: name line6 -> line7, C 5->14, D 12->11
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11

# Row 9: name extr1 -> line8, C 5->16, E 0->1
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("E9").Value = $true

# Row 10: name extr2 -> extr1, C 10->5, D 11->12, E 0->1
$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

# Row 11: name extr3 -> extr2, C 7->5, D 8->9
$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9

# Row 12: name extr4 -> extr3, C 9->10
$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10

# Row 13: name extr5 -> extr4, D 11->8
$ws.Range("B13").Value = "extr4"
$ws.Range("D13").Value = 8

# Row 14: name extr6 -> extr5, C 5->9, D 7->11, E 0->1
$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $true

# Row 15: name extr7 -> extr6, C (new)7, D(new)11, E(new)0
$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $false

# --- New rows 16 and 17 (copy formatting from row 15) ---
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true
